$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F, shifting the existing "District" column (F) to G.
$ws.Columns.Item(6).Insert()

# New header for the inserted "Address" column.
$ws.Cells.Item(2, 6).Value = "Address"

# Row 3 sub-header stays blank under "Address" (matches source: F3 has no text).
$ws.Cells.Item(3, 6).Value = ""

# Per-row school "Address" values, derived from the second line of column B
# (school name + taluk, with separating comma/space removed, trailing
# ", Raichur." dropped). Row 28 is left blank in the Address column because
# its original single-segment text ("G G U H S Fort Raichur") was already
# shifted into the District column (G28) by the column insert above.
$addresses = @{
    4  = "G H S KalmangiSindhnur"
    5  = "G H S GabburDevadurga"
    6  = "G H S BalaganurSindhanur"
    7  = "G H S TalekhanMaski"
    8  = "G U H S Police Colony"
    9  = "G H S Siyatalab"
    10 = "Govt. P U College(H S) GuruguntaLingasugur"
    11 = "G H S TurvihalSindhanur"
    12 = "Govt. Adarsha VidyalayaDeodurga"
    13 = "G U H S Arab Mohalla"
    14 = "G U H S Idapnur"
    15 = "G H S Jagir Venkatapur"
    16 = "G H S KasabalingasugurLingasugur"
    17 = "G H S Bijanageera"
    18 = "Govt. P U College (High School Section) Sirwar"
    19 = "G H S NagarahalLingasagur"
    20 = "G H S NagadadinniDeodurga"
    21 = "G H S Turkundona"
    22 = "G H S TimmapurSindhanur"
    23 = "U G H P S SannahosuruManvi"
    24 = "U G H P S Sagamkunta"
    25 = "G H S Yapaladinni"
    26 = "U G H P S KaradiguddaManvi"
    27 = "G H S Nilogal CrossLingasugur"
    29 = "G H S (B)JalahalliDevadurga"
    30 = "G H S Gillesugur"
    31 = "G H S KachapurLingsugur"
    32 = "G H S ChitapurLingasugur"
    33 = "G H S GonawarSindhanur"
    34 = "G H S YaradonLingasgur"
    35 = "Adarsha Vidyalaya Sindhanoor"
    36 = "G H S RajolliManvi"
    37 = "G H S SomanamaradiDevadurga"
    38 = "U G H R S MuddamuguddiManvi"
    39 = "U G H P S Urdu Manvi"
    40 = "G H S AmbhamathSindhanur"
}

foreach ($row in 4..40) {
    if ($addresses.ContainsKey($row)) {
        $ws.Cells.Item($row, 6).Value = $addresses[$row]
    }
}
